$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'291.07"
$ws.Range("E2").Value = "'1.40%"
$ws.Range("D3").Value = "'29.23"
$ws.Range("E3").Value = "'2.05%"
$ws.Range("D4").Value = "'5.284"
$ws.Range("E4").Value = "'4.92%"
$ws.Range("D5").Value = "'0.07140"
$ws.Range("E5").Value = "'7.25%"
$ws.Range("D6").Value = "'7.471"
$ws.Range("E6").Value = "'1.59%"
$ws.Range("D7").Value = "'3.570"
$ws.Range("E7").Value = "'5.39%"
$ws.Range("D8").Value = "'1.411"
$ws.Range("E8").Value = "'2.97%"
$ws.Range("D9").Value = "'0.9077"
$ws.Range("E9").Value = "'-3.62%"
$ws.Range("D10").Value = "'0.1624"
$ws.Range("E10").Value = "'4.16%"
$ws.Range("D11").Value = "'0.07697"
$ws.Range("E11").Value = "'15.90%"
$ws.Range("D12").Value = "'0.07725"
$ws.Range("E12").Value = "'1.10%"
$ws.Range("D13").Value = "'0.02921"
$ws.Range("E13").Value = "'-1.14%"
$ws.Range("D14").Value = "'0.09022"
$ws.Range("E14").Value = "'0.25%"
$ws.Range("D15").Value = "'0.001592"
$ws.Range("E15").Value = "'0.04%"
$ws.Range("D16").Value = "'0.0006517"
$ws.Range("E16").Value = "'0.57%"
$ws.Range("D17").Value = "'0.006198"
$ws.Range("E17").Value = "'-2.28%"
$ws.Range("D18").Value = "'3.482"
$ws.Range("E18").Value = "'1.04%"
$ws.Range("E19").Value = "'-1.18%"
$ws.Range("D20").Value = "'0.3243"
$ws.Range("E20").Value = "'0.84%"
$ws.Range("E21").Value = "'4.04%"
$ws.Range("D22").Value = "'3.980"
$ws.Range("E22").Value = "'-2.57%"
$ws.Range("E23").Value = "'2.91%"
$ws.Range("D24").Value = "'0.04511"
$ws.Range("E24").Value = "'0.43%"
$ws.Range("D25").Value = "'0.001203"
$ws.Range("E25").Value = "'2.05%"
$ws.Range("D26").Value = "'0.004203"
$ws.Range("E26").Value = "'-6.46%"
$ws.Range("D27").Value = "'0.0001159"
$ws.Range("E27").Value = "'-7.31%"
$ws.Range("D28").Value = "'0.0001673"
$ws.Range("E28").Value = "'3.37%"
$ws.Range("D40").Value = "'0.04404"
$ws.Range("E40").Value = "'4.74%"
$ws.Range("D41").Value = "'0.006816"
$ws.Range("E41").Value = "'0.71%"
$ws.Range("D42").Value = "'0.1263"
$ws.Range("E42").Value = "'0.50%"
$ws.Range("D43").Value = "'0.002069"
$ws.Range("E43").Value = "'2.43%"
$ws.Range("E44").Value = "'-4.45%"
$ws.Range("D45").Value = "'0.00005846"
$ws.Range("E45").Value = "'2.97%"
$ws.Range("D47").Value = "'0.01297"
$ws.Range("E47").Value = "'-0.74%"
